$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Info")

# New "Jenkins" column (E) on the Login Info sheet, header + flag value,
# styled to match the existing header / flag cells (D1 / D2).
$ws.Range("E1").Value = "Jenkins"
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4108

$ws.Range("E2").Value = 1
$ws.Range("E2").HorizontalAlignment = -4108

# Leave the active selection on D9, matching where the author's cursor was.
$ws.Range("D9").Select()
